{"js": "// Update the worksheet date (title paragraph) and all 100 two-digit\n// multiplication problems in the single table, in document order.\n// Some problem texts repeat (e.g. \"76\u00d725=\" appears twice) but map to\n// different replacements, so cells are addressed positionally\n// (row-major, matching the table's reading order) rather than via a\n// global find/replace.\n\nconst NEW_DATE = \"2023-07-24 Monday\";\n\n// 100 replacement values, row-major (20 rows x 5 columns), matching\n// the table's existing layout.\nconst NEW_CELLS_FLAT = [\"85\u00d736=\", \"89\u00d741=\", \"49\u00d799=\", \"55\u00d782=\", \"57\u00d745=\", \"21\u00d766=\", \"99\u00d770=\", \"54\u00d771=\", \"58\u00d718=\", \"72\u00d719=\", \"74\u00d755=\", \"10\u00d757=\", \"33\u00d767=\", \"67\u00d750=\", \"60\u00d741=\", \"61\u00d760=\", \"91\u00d798=\", \"70\u00d792=\", \"15\u00d786=\", \"78\u00d771=\", \"11\u00d743=\", \"24\u00d769=\", \"89\u00d750=\", \"83\u00d735=\", \"57\u00d761=\", \"95\u00d750=\", \"40\u00d749=\", \"61\u00d741=\", \"53\u00d738=\", \"37\u00d714=\", \"38\u00d748=\", \"33\u00d747=\", \"32\u00d743=\", \"19\u00d777=\", \"77\u00d721=\", \"74\u00d724=\", \"49\u00d720=\", \"22\u00d731=\", \"47\u00d797=\", \"48\u00d724=\", \"95\u00d796=\", \"24\u00d724=\", \"21\u00d724=\", \"19\u00d729=\", \"16\u00d797=\", \"65\u00d729=\", \"54\u00d756=\", \"73\u00d783=\", \"35\u00d785=\", \"28\u00d7100=\", \"25\u00d797=\", \"45\u00d739=\", \"83\u00d765=\", \"85\u00d722=\", \"73\u00d740=\", \"29\u00d732=\", \"87\u00d757=\", \"99\u00d795=\", \"75\u00d745=\", \"39\u00d730=\", \"35\u00d757=\", \"97\u00d752=\", \"39\u00d776=\", \"74\u00d763=\", \"46\u00d7100=\", \"76\u00d719=\", \"78\u00d742=\", \"48\u00d755=\", \"16\u00d768=\", \"73\u00d795=\", \"55\u00d724=\", \"76\u00d795=\", \"18\u00d770=\", \"79\u00d720=\", \"51\u00d761=\", \"83\u00d711=\", \"59\u00d775=\", \"58\u00d756=\", \"84\u00d799=\", \"57\u00d777=\", \"32\u00d753=\", \"56\u00d745=\", \"48\u00d798=\", \"94\u00d728=\", \"66\u00d770=\", \"19\u00d797=\", \"25\u00d740=\", \"61\u00d799=\", \"38\u00d795=\", \"54\u00d732=\", \"79\u00d764=\", \"58\u00d724=\", \"31\u00d739=\", \"33\u00d791=\", \"82\u00d784=\", \"28\u00d780=\", \"88\u00d712=\", \"34\u00d712=\", \"71\u00d770=\", \"28\u00d744=\"];\n\n// --- Title paragraph (the date line) -------------------------------\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\ntitlePara.insertText(NEW_DATE, \"Replace\");\n\n// --- Table of multiplication problems -------------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst columnCount = 5;\nlet k = 0;\nfor (let r = 0; r < rows.items.length; r++) {\n  const cells = rows.items[r].cells;\n  cells.load(\"items/value\");\n  await context.sync();\n  for (let c = 0; c < cells.items.length; c++) {\n    cells.items[c].value = NEW_CELLS_FLAT[k];\n    k++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date (first paragraph) and all 100 two-digit\n# multiplication problems in the single table, in document order.\n# Some problem texts repeat (e.g. \"76x25=\" appears twice) but map to\n# different replacements, so cells are addressed positionally by\n# (row, column) rather than via a global Find/Replace.\n\n$d = $word.ActiveDocument\n\n$NEW_DATE = \"2023-07-24 Monday\"\n\n$NEW_VALUES = @(\n    @(\"85\u00d736=\", \"89\u00d741=\", \"49\u00d799=\", \"55\u00d782=\", \"57\u00d745=\"),\n    @(\"21\u00d766=\", \"99\u00d770=\", \"54\u00d771=\", \"58\u00d718=\", \"72\u00d719=\"),\n    @(\"74\u00d755=\", \"10\u00d757=\", \"33\u00d767=\", \"67\u00d750=\", \"60\u00d741=\"),\n    @(\"61\u00d760=\", \"91\u00d798=\", \"70\u00d792=\", \"15\u00d786=\", \"78\u00d771=\"),\n    @(\"11\u00d743=\", \"24\u00d769=\", \"89\u00d750=\", \"83\u00d735=\", \"57\u00d761=\"),\n    @(\"95\u00d750=\", \"40\u00d749=\", \"61\u00d741=\", \"53\u00d738=\", \"37\u00d714=\"),\n    @(\"38\u00d748=\", \"33\u00d747=\", \"32\u00d743=\", \"19\u00d777=\", \"77\u00d721=\"),\n    @(\"74\u00d724=\", \"49\u00d720=\", \"22\u00d731=\", \"47\u00d797=\", \"48\u00d724=\"),\n    @(\"95\u00d796=\", \"24\u00d724=\", \"21\u00d724=\", \"19\u00d729=\", \"16\u00d797=\"),\n    @(\"65\u00d729=\", \"54\u00d756=\", \"73\u00d783=\", \"35\u00d785=\", \"28\u00d7100=\"),\n    @(\"25\u00d797=\", \"45\u00d739=\", \"83\u00d765=\", \"85\u00d722=\", \"73\u00d740=\"),\n    @(\"29\u00d732=\", \"87\u00d757=\", \"99\u00d795=\", \"75\u00d745=\", \"39\u00d730=\"),\n    @(\"35\u00d757=\", \"97\u00d752=\", \"39\u00d776=\", \"74\u00d763=\", \"46\u00d7100=\"),\n    @(\"76\u00d719=\", \"78\u00d742=\", \"48\u00d755=\", \"16\u00d768=\", \"73\u00d795=\"),\n    @(\"55\u00d724=\", \"76\u00d795=\", \"18\u00d770=\", \"79\u00d720=\", \"51\u00d761=\"),\n    @(\"83\u00d711=\", \"59\u00d775=\", \"58\u00d756=\", \"84\u00d799=\", \"57\u00d777=\"),\n    @(\"32\u00d753=\", \"56\u00d745=\", \"48\u00d798=\", \"94\u00d728=\", \"66\u00d770=\"),\n    @(\"19\u00d797=\", \"25\u00d740=\", \"61\u00d799=\", \"38\u00d795=\", \"54\u00d732=\"),\n    @(\"79\u00d764=\", \"58\u00d724=\", \"31\u00d739=\", \"33\u00d791=\", \"82\u00d784=\"),\n    @(\"28\u00d780=\", \"88\u00d712=\", \"34\u00d712=\", \"71\u00d770=\", \"28\u00d744=\")\n)\n\n# --- Title paragraph (the date line) -------------------------------\n$d.Paragraphs.Item(1).Range.Text = $NEW_DATE\n\n# --- Table of multiplication problems -------------------------------\n$table = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $NEW_VALUES.Count; $r++) {\n    $rowValues = $NEW_VALUES[$r - 1]\n    for ($c = 1; $c -le $rowValues.Count; $c++) {\n        $table.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
